$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6..24 down to 7..25
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly price record
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44952
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100103
$ws.Cells.Item(6, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(6, 9).Value = 100103002
$ws.Cells.Item(6, 10).Value = "Ciruela"
$ws.Cells.Item(6, 11).Value = "Larry Ann"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 21000
$ws.Cells.Item(6, 16).Value = 20500
$ws.Cells.Item(6, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(6, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(6, 19).Value = 1139
$ws.Cells.Item(6, 20).Value = 18
